$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows to append (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$newRows = @(
    @(44308, 5, 49, 191.7883283103057),
    @(44309, 12, 51, 199.6164233433794),
    @(44310, 10, 50, 195.7023758268425),
    @(44311, 17, 60, 234.842850992211),
    @(44312, 5, 61, 238.7568985087479)
)

$startRow = 234
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $newRows[$i]

    $ws.Range($ws.Cells.Item($r - 1, 1), $ws.Cells.Item($r - 1, 1)).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $rowData[0]

    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
}
